$d = $word.ActiveDocument

# Find the end of the "LOQ4073: Química Geral II (Requisito fraco)" paragraph
# (the last requirement line that should survive) and extend past its own
# paragraph mark so deletion starts at the very next paragraph.
$startRange = $d.Content
$startRange.Find.Execute("LOQ4073: Química Geral II (Requisito fraco)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startRange.MoveEnd(1, 1)
$startPos = $startRange.End

# Find the end of the copyright/footer paragraph that should be removed,
# and extend past its own paragraph mark too, so the trailing blank
# paragraph + page-break paragraph remain untouched.
$endRange = $d.Content
$endRange.Find.Execute("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endRange.MoveEnd(1, 1)
$endPos = $endRange.End

# Delete the blank paragraph, the "Ver no Jupiter..." paragraph, and the
# copyright/footer paragraph in one shot.
$deleteRange = $d.Range($startPos, $endPos)
$deleteRange.Delete()
